$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "NCT(2.785147743511865, 1.2057023122632777, -0.007301044594589297, 2.3184507413268456)"
$ws.Range("C2").Value = "JSU(-1.3658504100454936, 1.17601274990059, 2.319125548339838, 4.550893059250333)"
$ws.Range("D2").Value = "NCT(2.519138760464809, 1.2966275864424683, -0.4729688912212425, 2.4628599625295866)"
$ws.Range("E2").Value = "NIG(0.9868407308973548, 0.797782698808083, 4.001152984074338, 4.2282621873250275)"
